$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header label: "Asset Sub Class ID" -> "Asset Sub Class"
$ws.Range("A1").Value = "Asset Sub Class"

# Fix A2 value: numeric 2 -> text "Buildings"
$ws.Range("A2").Value = "Buildings"

# Update selection (logo export view state) from D10 to H8
$ws.Range("H8").Select()
